# Slide 2 ("Project overview"): bump the body text placeholder's font
# size to 24pt across all its runs/paragraphs.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

$shp = $null
foreach ($candidate in $s.Shapes) {
    if ($candidate.Name -eq "Content Placeholder 2") {
        $shp = $candidate
        break
    }
}
if ($shp -eq $null) {
    $shp = $s.Shapes.Item(2)
}

$tr = $shp.TextFrame.TextRange
$tr.Font.Size = 24
